$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Control system")
$ws.Range("AM2").Value = 1
$ws.Range("AM2").HorizontalAlignment = -4108
$ws.Range("AM2").Borders.LineStyle = 1
